# Daily attendance processing - 2026-01-30 09:50:53
#
# The "Recorded By" column (G) lists the session recorder(s) as a
# comma-separated string. Rows that were recorded by both the automated
# "System" process and the user "dnasr281@gmail.com" currently read
# "System, dnasr281@gmail.com"; flip the order to
# "dnasr281@gmail.com, System" throughout the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

[void]$ws.Cells.Replace("System, dnasr281@gmail.com", "dnasr281@gmail.com, System", 1)
